$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 25.23990433333333
$ws.Range("H2").Value = 75.719713
$ws.Range("I2").Value = 0.05173702626903214
$ws.Range("J2").Value = 0.05173702626903214
$ws.Range("M2").Value = 0.05031533333333333
$ws.Range("N2").Value = 0.150946
$ws.Range("O2").Value = 0.005485022167780355
$ws.Range("P2").Value = 0.005485022167780356
$ws.Range("Q2").Value = 1.269954199833111
$ws.Range("R2").Value = 11.429587798498
$ws.Range("S2").Value = 0.0002837787359806758
$ws.Range("T2").Value = 0.0002837787359806759

# Row 3
$ws.Range("G3").Value = 25.23990433333333
$ws.Range("H3").Value = 75.719713
$ws.Range("I3").Value = 0.05173702626903214
$ws.Range("J3").Value = 0.05173702626903214
$ws.Range("O3").Value = 0.377950825718477
$ws.Range("P3").Value = 0.377950825718477
$ws.Range("Q3").Value = 87.50743821438512
$ws.Range("R3").Value = 787.566943929466
$ws.Range("S3").Value = 0.01955405179859923
$ws.Range("T3").Value = 0.01955405179859923

# Row 4
$ws.Range("G4").Value = 25.23990433333333
$ws.Range("H4").Value = 75.719713
$ws.Range("I4").Value = 0.05173702626903214
$ws.Range("J4").Value = 0.05173702626903214
$ws.Range("M4").Value = 5.655880666666666
$ws.Range("N4").Value = 16.967642
$ws.Range("O4").Value = 0.6165641521137426
$ws.Range("P4").Value = 0.6165641521137426
$ws.Range("Q4").Value = 142.7538869474162
$ws.Range("R4").Value = 1284.784982526746
$ws.Range("S4").Value = 0.03189919573445223
$ws.Range("T4").Value = 0.03189919573445223

# Row 5
$ws.Range("I5").Value = 0.8454897015965644
$ws.Range("J5").Value = 0.8454897015965646
$ws.Range("M5").Value = 0.05031533333333333
$ws.Range("N5").Value = 0.150946
$ws.Range("O5").Value = 0.005485022167780355
$ws.Range("P5").Value = 0.005485022167780356
$ws.Range("Q5").Value = 20.75367053132889
$ws.Range("R5").Value = 186.78303478196
$ws.Range("S5").Value = 0.004637529755887153
$ws.Range("T5").Value = 0.004637529755887155

# Row 6
$ws.Range("I6").Value = 0.8454897015965644
$ws.Range("J6").Value = 0.8454897015965646
$ws.Range("O6").Value = 0.377950825718477
$ws.Range("P6").Value = 0.377950825718477
$ws.Range("S6").Value = 0.3195535308548903
$ws.Range("T6").Value = 0.3195535308548904

# Row 7
$ws.Range("I7").Value = 0.8454897015965644
$ws.Range("J7").Value = 0.8454897015965646
$ws.Range("M7").Value = 5.655880666666666
$ws.Range("N7").Value = 16.967642
$ws.Range("O7").Value = 0.6165641521137426
$ws.Range("P7").Value = 0.6165641521137426
$ws.Range("Q7").Value = 2332.892900517657
$ws.Range("R7").Value = 20996.03610465892
$ws.Range("S7").Value = 0.521298640985787
$ws.Range("T7").Value = 0.521298640985787

# Row 8
$ws.Range("G8").Value = 50.137933
$ws.Range("H8").Value = 150.413799
$ws.Range("I8").Value = 0.1027732721344034
$ws.Range("J8").Value = 0.1027732721344034
$ws.Range("M8").Value = 0.05031533333333333
$ws.Range("N8").Value = 0.150946
$ws.Range("O8").Value = 0.005485022167780355
$ws.Range("P8").Value = 0.005485022167780356
$ws.Range("Q8").Value = 2.522706811539333
$ws.Range("R8").Value = 22.704361303854
$ws.Range("S8").Value = 0.0005637136759125255
$ws.Range("T8").Value = 0.0005637136759125256

# Row 9
$ws.Range("G9").Value = 50.137933
$ws.Range("H9").Value = 150.413799
$ws.Range("I9").Value = 0.1027732721344034
$ws.Range("J9").Value = 0.1027732721344034
$ws.Range("O9").Value = 0.377950825718477
$ws.Range("P9").Value = 0.377950825718477
$ws.Range("Q9").Value = 173.8295841478353
$ws.Range("R9").Value = 1564.466257330518
$ws.Range("S9").Value = 0.0388432430649875
$ws.Range("T9").Value = 0.0388432430649875

# Row 10
$ws.Range("G10").Value = 50.137933
$ws.Range("H10").Value = 150.413799
$ws.Range("I10").Value = 0.1027732721344034
$ws.Range("J10").Value = 0.1027732721344034
$ws.Range("M10").Value = 5.655880666666666
$ws.Range("N10").Value = 16.967642
$ws.Range("O10").Value = 0.6165641521137426
$ws.Range("P10").Value = 0.6165641521137426
$ws.Range("Q10").Value = 283.5741659213286
$ws.Range("R10").Value = 2552.167493291957
$ws.Range("S10").Value = 0.06336631539350333
$ws.Range("T10").Value = 0.06336631539350333

Write-Output "Updated TPM values"